$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) updates.
# Some new prices look like plain numbers (e.g. "7.40"); the source data
# stores the Price/Volume columns as text (e.g. to keep trailing zeros and
# thousands-dot formatting like "61.734.53"), so format those cells as Text
# first to prevent Excel from auto-converting them to numeric values.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = "61.734.53"
$ws.Range("D3").Value = "3.412.63"
$ws.Range("D4").Value = "0.999"
$ws.Range("D5").Value = "577.33"
$ws.Range("D6").Value = "144.34"
$ws.Range("D9").Value = "7.63"
$ws.Range("D11").Value = "0.384"
$ws.Range("D12").Value = "3.994.16"
$ws.Range("D14").Value = "27.97"
$ws.Range("D15").Value = "3.411.45"
$ws.Range("D16").Value = "0.0000169"
$ws.Range("D17").Value = "61.757.18"
$ws.Range("D19").Value = "13.83"
$ws.Range("D21").Value = "388.71"
$ws.Range("D22").Value = "74.20"
$ws.Range("D23").Value = "0.549"
$ws.Range("D25").Value = "0.0000114"
$ws.Range("D26").Value = "0.188"
$ws.Range("D27").Value = "0.999"
$ws.Range("D28").Value = "7.40"
$ws.Range("D33").Value = "23.44"
$ws.Range("D34").Value = "5.17"
$ws.Range("D35").Value = "6.93"
$ws.Range("D36").Value = "168.46"
$ws.Range("D37").Value = "3.443.59"
$ws.Range("D38").Value = "1.47"
$ws.Range("D39").Value = "28.46"
$ws.Range("D41").Value = "0.786"
$ws.Range("D44").Value = "1.17"
$ws.Range("D45").Value = "2.504.02"
$ws.Range("D46").Value = "22.68"
$ws.Range("D47").Value = "0.999"
$ws.Range("D48").Value = "6.60"
$ws.Range("D49").Value = "0.0262"
$ws.Range("D50").Value = "2.09"

# Volume(1h) (column E) updates
$ws.Range("E2").Value = "  +1.24%  "
$ws.Range("E3").Value = "  +0.84%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("E5").Value = "  +1.08%  "
$ws.Range("E6").Value = "  +1.77%  "
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("E9").Value = "  -0.05%  "
$ws.Range("E10").Value = "  +0.12%  "
$ws.Range("E12").Value = "  +0.77%  "
$ws.Range("E13").Value = "  -0.52%  "
$ws.Range("E14").Value = "  +0.23%  "
$ws.Range("E15").Value = "  +0.51%  "
$ws.Range("E16").Value = "  -1.06%  "
$ws.Range("E17").Value = "  +1.11%  "
$ws.Range("E19").Value = "  +1.65%  "
$ws.Range("E20").Value = "  +2.75%  "
$ws.Range("E21").Value = "  +1.36%  "
$ws.Range("E22").Value = "  -1.32%  "
$ws.Range("E23").Value = "  -0.64%  "
$ws.Range("E24").Value = "  +0.06%  "
$ws.Range("E25").Value = "  -1.41%  "
$ws.Range("E26").Value = "  +2.89%  "
$ws.Range("E27").Value = "  -0.04%  "
$ws.Range("E28").Value = "  +1.94%  "
$ws.Range("E29").Value = "  +0.40%  "
$ws.Range("E30").Value = "  +0.02%  "
$ws.Range("E31").Value = "  +1.15%  "
$ws.Range("E32").Value = "  +0.00%  "
$ws.Range("E33").Value = "  +1.01%  "
$ws.Range("E34").Value = "  +3.93%  "
$ws.Range("E35").Value = "  -0.08%  "
$ws.Range("E36").Value = "  +1.21%  "
$ws.Range("E37").Value = "  +0.80%  "
$ws.Range("E38").Value = "  +0.24%  "
$ws.Range("E39").Value = "  +6.30%  "
$ws.Range("E40").Value = "  -1.58%  "
$ws.Range("E41").Value = "  +0.87%  "
$ws.Range("E42").Value = "  +1.57%  "
$ws.Range("E43").Value = "  +0.63%  "
$ws.Range("E44").Value = "  +4.07%  "
$ws.Range("E45").Value = "  +2.35%  "
$ws.Range("E46").Value = "  -0.92%  "
$ws.Range("E47").Value = "  -0.06%  "
$ws.Range("E48").Value = "  -1.20%  "
$ws.Range("E49").Value = "  -0.35%  "
$ws.Range("E50").Value = "  -3.15%  "
$ws.Range("E51").Value = "  -0.58%  "
